$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
